# Apply the changes described by the commit:
#  - active tab moves from sheet index 1 ("test") to sheet index 2 ("test2")
#  - tabSelected moves from "test" to "test2" (handled by activating/selecting "test2")
#  - column width tweaks on all three sheets
#  - new cells E1/F1/G1 on sheet "test2", introducing a new shared string "newname"

$wb = $excel.ActiveWorkbook

$wsDemo  = $wb.Worksheets.Item("demo")
$wsTest  = $wb.Worksheets.Item("test")
$wsTest2 = $wb.Worksheets.Item("test2")

# Column width tweaks (values chosen so that the stored OOXML width lands as
# close as possible to the target width used by the authoring application).
$wsDemo.Columns.Item(1).ColumnWidth  = 25.714285714285715
$wsTest.Columns.Item(1).ColumnWidth  = 10.571428571428571
$wsTest2.Columns.Item(1).ColumnWidth = 9.857142857142858

# New data on sheet "test2", row 1
$wsTest2.Range("E1").Value = "^"
$wsTest2.Range("F1").Value = "name"
$wsTest2.Range("G1").Value = "newname"

# Make "test2" the active / selected tab (was "test" before).
$wsTest2.Select()
$wsTest2.Activate()
